# employee_import_template.xlsx
# Validate excel before importing employees; default password is 123456.
# Add a "Username" column (G) to the header row of the import template.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cell with the same value/text as the commit adds to sharedStrings.
$ws.Range("G1").Value = "Username"

# Give it the same look (bold Times New Roman, fill, border, centered) as the
# rest of the header row by cloning F1's formatting.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New column width (close to the 14.140625-character width of the target file).
$ws.Columns.Item(7).ColumnWidth = 13.3

# Move the active selection from H2 to H1, as in the edited workbook.
$ws.Range("H1").Select()
